$d = $word.ActiveDocument

$result0 = $d.Content.Find.Execute("Mechanical Testing", $true, $false, $false, $false, $false, $true, 1, $false, "Mechanical Behavior of Materials", 2)
Write-Host "Replace 0: $result0"

$result1 = $d.Content.Find.Execute("Créditos-aula: 2", $true, $false, $false, $false, $false, $true, 1, $false, "Créditos-aula: 4", 2)
Write-Host "Replace 1: $result1"

$result2 = $d.Content.Find.Execute("Carga horária: 30 h", $true, $false, $false, $false, $false, $true, 1, $false, "Carga horária: 60 h", 2)
Write-Host "Replace 2: $result2"

$result3 = $d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2)
Write-Host "Replace 3: $result3"

$result4 = $d.Content.Find.Execute("Descrever os ensaios mecânicos usuais para determinar propriedades de rigidez elástica, resistência e ductilidade dos materiais metálicos, cerâmicos e poliméricos. Descrever os principais métodos experimentais e técnicas de análise envolvidas. Habilitar os alunos ao tratamento e análise dos resultados experimentais obtidos.", $true, $false, $false, $false, $false, $true, 1, $false, "Descrever os ensaios mecânicos usuais para determinar propriedades de rigidez elástica, resistência e ductilidade dosmateriais metálicos, cerâmicos e poliméricos. Descrever os principais métodos experimentais e técnicas de análiseenvolvidas. Habilitar os alunos ao tratamento e análise dos resultados experimentais obtidos pois a disciplina possui cunho teórico/prático com idas aos laboratórios para sedimentar a teoria vista em aulas expositivas.", 2)
Write-Host "Replace 4: $result4"

$result5 = $d.Content.Find.Execute("The discipline is aimed at describing the usual mechanical testing methods to determine the properties of elastic stiffness, strength and ductility of metallic, ceramic and polymeric materials. The main experimental and data analysis techniques are presented, enabling students to assess the obtained experimental results.", $true, $false, $false, $false, $false, $true, 1, $false, "The discipline is aimed at describing the usual mechanical testing methods to determine the properties of elastic stiffness, strength and ductility of metallic, ceramic and polymeric materials. Enable students to process and analyze the experimental results obtained, as the discipline has a theoretical/experimental nature with visits to laboratories to solidify the theory learned in lectures.", 2)
Write-Host "Replace 5: $result5"

$result6 = $d.Content.Find.Execute("Importância dos ensaios mecânicos; Introdução à extensometria; Medidas de dureza; Ensaio de tração; Ensaio de Flexão; Ensaios de compressão, torção, e resistência à tração por compressão diametral; Fratura frágil e o ensaio de impacto; Introdução aos ensaios de fadiga e de fluência.", $true, $false, $false, $false, $false, $true, 1, $false, "Importância dos ensaios mecânicos; Introdução à extensometria; Medidas de dureza; Ensaio de tração; Ensaio de Flexão; Ensaios de compressão, Ensaios de torção, e resistência à tração por compressão diametral; Fratura dúctil e frágil e o ensaio de impacto Charpy e Izod.", 2)
Write-Host "Replace 6: $result6"

$result7 = $d.Content.Find.Execute("Significance of mechanical testing; Introduction to strain gage fundamentals; Hardness testing; Tensile testing; Bending testing; Standard methods for compression, torsion and splitting tensile strength determination; Brittle fracture and the impact test; Introduction to fatigue and creep testing methods.", $true, $false, $false, $false, $false, $true, 1, $false, "Significance of mechanical testing; Introduction to strain gage fundamentals; Hardness testing; Tensile testing; Bending testing; Standard methods for compression, torsion and splitting tensile strength determination; Ductile and Brittle fracture and the charpy and izod impact testing methods.", 2)
Write-Host "Replace 7: $result7"

$result8 = $d.Content.Find.Execute("1.Importância dos ensaios e do estabelecimento de normas para sua execução; classificação dos ensaios. 2.Introdução à extensometria: análise experimental de tensões e deformações. 3.Medidas de dureza: Conceitos e tipos de ensaio. Ensaio de dureza por penetração: principais escalas e equipamentos. Ensaio de microdureza. Dureza ao choque. 4.O ensaio de tração. Parâmetros de ensaio. Propriedades obtidas no ensaio. Curvas tensão-deformação de engenharia e real. 5.Ensaios de flexão a três e quatro pontos: determinação do módulo de ruptura e módulo de Young, Fontes de erros no ensaio de flexão e como minimizá-los. Estatística de Weibull; 6.Ensaio de compressão: técnica experimental e cuidados necessários. Compressão de materiais dúcteis e frágeis. Resistência à tração por compressão diametral. Ensaio de torção: Determinação do módulo de elasticidade transversal. 7.Ensaio de impacto. Fatores que concorrem para a fratura frágil dos materiais. A transição dúctil-frágil. Métodos Charpy e Izod. Ensaio de impacto instrumentado. 8. Introdução aos ensaios de fadiga e fluência. Métodos de ensaio de fadiga de alto-ciclo e baixo-ciclo. Crescimento da trinca por fadiga. Critérios de projeto em fadiga. Métodos experimentais para o ensaio de fluência. Formas de representação dos resultados. Taxa de fluência estacionária e determinação da energia de ativação para fluência.", $true, $false, $false, $false, $false, $true, 1, $false, "Programa1.Importância dos ensaios e do estabelecimento de normas para sua execução; classificação dos ensaios. 2.Introdução à extensometria: análise experimental de tensões e deformações. 3.Medidas de dureza: Conceitos e tipos de ensaio. Ensaio de dureza por penetração: principais escalas e equipamentos. Ensaio de microdureza. Dureza ao choque. 4.O ensaio de tração. Parâmetros de ensaio. Propriedades obtidas no ensaio. Curvas tensão-deformação de engenharia e real. 5. Ensaios de flexão a três e quatro pontos: determinação do módulo de ruptura e módulo de Young, Fontes de erros no ensaio de flexão e como minimizá-los. Estatística de Weibull; 6. Ensaio de compressão: técnica experimental e cuidados necessários. Compressão de materiais dúcteis e frágeis. Resistência à tração por compressão diametral. Ensaio de torção: Determinação do módulo de elasticidade transversal. 7.Ensaio de impacto. Fatores que concorrem para a fratura frágil dos materiais. A transição dúctil-frágil. Métodos Charpy e Izod. Ensaio de impacto instrumentado. 8. Métodos experimentais para o ensaio de fluência. Formas de representação dos resultados. Taxa de fluência estacionária e determinação da energia de ativação para fluência.", 2)
Write-Host "Replace 8: $result8"

$result9 = $d.Content.Find.Execute("1. Significance of mechanical testing and the standard testing procedures; Classes of testing methods. 2. Introduction to strain gage fundamentals: experimental stress and strain analysis. 3. Hardness measurements: types of hardness tests. Indentation hardness: usual scales and testing machines. Microhardness test. Rebound hardness test. 4. Tensile test: test parameters and properties determination. Engineering and true stress-strain curves. 5. Three- and Four-point bending tests: determination of Young’s modulus and modulus of rupture. Error sources and minimization. Weibull statistics. 6. Compression test: experimental method. Compression of brittle and ductile materials. The splitting tensile strength determination. Torsion test: determination of the shear modulus. 7.Impact testing. Factors that contribute to brittle fracture. The ductile-brittle transition. Charpy and Izod testing methods. Instrumented impact test. 8. Introduction to fatigue and creep testing. High- and Low-cycle fatigue. Fatigue crack propagation. Fatigue design criteria. Experimental methods for creep testing. Representing creep tests results. The steady-state creep rate and the activation energy for creep.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Significance of mechanical testing and the standard testing procedures; Classes of testing methods. 2. Introduction to strain gage fundamentals: experimental stress and strain analysis. 3. Hardness measurements: types of hardness tests. Indentation hardness: usual scales and testing machines. Microhardness test. Rebound hardness test. 4. Tensile test: test parameters and properties determination. Engineering and true stress-strain curves. 5. Three- and Four-point bending tests: determination of Young’s modulus and modulus of rupture. Error sources and minimization. Weibull statistics. 6. Compression test: experimental method. Compression of brittle and ductile materials. The splitting tensile strength determination. Torsion test: determination of the shear modulus. 7. Impact testing. Factors that contribute to brittle fracture. The ductile-brittle transition. Charpy and Izod testing methods. Instrumented impact test. 8. Experimental methods for creep testing. Representing creep tests results. The stead-state creep rate and the activation energy for creep.", 2)
Write-Host "Replace 9: $result9"
